$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers (U1:Y1) ---
$ws.Range("U1").Value = "Planet"
$ws.Range("V1").Value = "Resultant x"
$ws.Range("W1").Value = "Resultant y"
$ws.Range("X1").Value = "Acceleration x"
$ws.Range("Y1").Value = "Acceleration y"

# --- Column widths for the new columns (U..Y) ---
# ColumnWidth (character units) is stored internally as (value + 5/6) by this
# engine, so back the requested value off by 5/6 to land on the target.
$ws.Columns.Item(21).ColumnWidth = 12.053385416666666   # U -> 12.88671875
$ws.Range("V1:W1").EntireColumn.ColumnWidth = 14.166666666666666  # V,W -> 15
$ws.Columns.Item(24).ColumnWidth = 16.166666666666668   # X -> 17
$ws.Columns.Item(25).ColumnWidth = 13.608072916666666   # Y -> 14.44140625

# --- Row 2 (planet 1) ---
$ws.Range("U2").Value = 1
$ws.Range("V2").Formula = "=SUM(N2, N3, N4, N5)"
$ws.Range("W2").Formula = "=SUM(O2, O3, O4, O5)"
$ws.Range("X2").Formula = "=V2/B2"
$ws.Range("Y2").Formula = "=W2/B2"

# --- Row 3 (planet 2) ---
$ws.Range("U3").Value = 2
$ws.Range("V3").Formula = "=SUM(N6,N7,N8,P2)"
$ws.Range("W3").Formula = "=SUM(O8,O7,O6,Q2)"

# --- Row 4 (planet 3) ---
$ws.Range("U4").Value = 3
$ws.Range("V4").Formula = "=SUM(N9,N10,P3,P6)"
$ws.Range("W4").Formula = "=SUM(O9,O10,Q3,Q6)"

# --- Row 5 (planet 4) ---
$ws.Range("U5").Value = 4
$ws.Range("V5").Formula = "=SUM(N11,P4,P7,P9)"
$ws.Range("W5").Formula = "=SUM(O11,Q9,Q7,Q4)"

# --- Row 6 (planet 5) ---
$ws.Range("U6").Value = 5
$ws.Range("V6").Formula = "=SUM(P11,P10,P8,P5)"
$ws.Range("W6").Formula = "=SUM(Q11,Q10,Q8,Q5)"

# X3:X6 / Y3:Y6 share one formula each across the range (matches the
# t="shared" formula groups introduced by the diff).
$ws.Range("X3:X6").Formula = "=V3/B3"
$ws.Range("Y3:Y6").Formula = "=W3/B3"

# --- Sheet view: scroll + selection moved to the new columns ---
$ws.Range("X5:Y5").Select()

# --- Page setup: portrait / letter-ish paper added by the edit ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
